# Trade #59 closed at 2026-02-17 08:48:06 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.62
$summary.Range("B4").Value = -0.38
$summary.Range("B5").Value = -0.13
$summary.Range("B6").Value = 59
$summary.Range("B7").Value = 23
$summary.Range("B9").Value = 38.98

# ---- Strategy Status sheet (MarketMaking row) ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.62
$status.Range("D4").Value = 59
$status.Range("E4").Value = -0.38
$status.Range("F4").Value = -0.38
$status.Range("G4").Value = 38.98

# ---- New trade row data (shared by All Trades + MarketMaking sheets) ----
function Add-TradeRow($ws, [int]$row) {
    $ws.Cells.Item($row, 1).Value = 59

    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.Value = "'2026-02-17"
    $dateCell.Style = "Normal"

    $timeCell = $ws.Cells.Item($row, 3)
    $timeCell.Value = "'08:48:00"
    $timeCell.Style = "Normal"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.67
    $ws.Cells.Item($row, 7).Value = 0.6899999999999999
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 2.9851
    $ws.Cells.Item($row, 10).Value = 0.02
    $ws.Cells.Item($row, 11).Value = 99.62
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 60

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 60
